$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo "Indi Project" -> "Indie Project" in the existing log entry (row 37)
$ws.Range("D37").Value = "Indie Project: consider how logic class will work - is it really a special dao?; Drafted generic dao"

# Insert three new rows after row 40 (the old empty placeholder row), pushing
# all later entries down by three rows. Excel copies the formatting of the
# row above (row 40: A s=1 date-column style, D s=7 empty time-style) into
# the newly inserted rows.
$ws.Rows("41:43").Insert()

# Rows 41 and 42 become plain blank spacer rows (only the date-column style
# on A carries over); drop the inherited D formatting entirely.
$ws.Range("D41").Clear()
$ws.Range("D42").Clear()

# Row 40 (previously the empty time-format placeholder) now holds the new
# log entry; row 43 becomes the new empty placeholder (keeps the style that
# used to live on D40).
$ws.Range("D40").Clear()
$ws.Range("D40").Value = "11:45 - x"

# Reflect the saved selection/scroll position from the authored workbook.
$ws.Range("D41").Select()
